$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 2466.1667
$ws.Range("I20").Value = 2466.1667
$ws.Range("K20").Value = 2466.1667
$ws.Range("M20").Value = -2236.1667
$ws.Range("H35").Value = 2466.1667
$ws.Range("I35").Value = 2466.1667
$ws.Range("K35").Value = 2466.1667
$ws.Range("M35").Value = -2087.1667
$ws.Range("H74").Value = 6983.857
$ws.Range("I74").Value = 3830.9092
$ws.Range("K74").Value = 3830.9092
$ws.Range("M74").Value = -2894.9092
$ws.Range("H77").Value = 6983.857
$ws.Range("I77").Value = 3830.9092
$ws.Range("K77").Value = 19154.546
$ws.Range("M77").Value = -14474.546
$ws.Range("H113").Value = 2955.5
$ws.Range("I113").Value = 3046.6
$ws.Range("K113").Value = 3046.6
$ws.Range("M113").Value = 207.4000000000001
$ws.Range("H132").Value = 49696.594
$ws.Range("I132").Value = 54353.066
$ws.Range("K132").Value = 163059.198
$ws.Range("M132").Value = -160529.198
$ws.Range("H141").Value = 1563.3334
$ws.Range("I141").Value = 1256
$ws.Range("J141").Value = 3100
$ws.Range("K141").Value = 3768
$ws.Range("L141").Value = 9300
$ws.Range("M141").Value = 1412
$ws.Range("N141").Value = -19660

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 1786.8572
$ws.Range("I30").Value = 1211.6
$ws.Range("K30").Value = 1211.6
$ws.Range("M30").Value = -1061.6
$ws.Range("H32").Value = 5955583.5
$ws.Range("I32").Value = 6413476.5
$ws.Range("K32").Value = 6413476.5
$ws.Range("M32").Value = -6413189.5
$ws.Range("H92").Value = 84958.336
$ws.Range("I92").Value = 70000
$ws.Range("J92").Value = 87950
$ws.Range("K92").Value = 70000
$ws.Range("L92").Value = 87950
$ws.Range("M92").Value = -67504
$ws.Range("N92").Value = -92942
$ws.Range("H105").Value = 90369
$ws.Range("J105").Value = 90369
$ws.Range("L105").Value = 90369
$ws.Range("N105").Value = -97357
$ws.Range("H132").Value = 643059.2
$ws.Range("I132").Value = 848182.1
$ws.Range("K132").Value = 2544546.3
$ws.Range("M132").Value = -2542016.3

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1633
$ws.Range("I86").Value = 1799.5
$ws.Range("J86").Value = 1300
$ws.Range("K86").Value = 1799.5
$ws.Range("L86").Value = 1300
$ws.Range("M86").Value = -676.5
$ws.Range("N86").Value = -3546
$ws.Range("H89").Value = 1633
$ws.Range("I89").Value = 1799.5
$ws.Range("J89").Value = 1300
$ws.Range("K89").Value = 8997.5
$ws.Range("L89").Value = 6500
$ws.Range("M89").Value = -3381.5
$ws.Range("N89").Value = -17732
$ws.Range("H94").Value = 943.5833
$ws.Range("I94").Value = 871.75
$ws.Range("J94").Value = 1087.25
$ws.Range("K94").Value = 871.75
$ws.Range("L94").Value = 1087.25
$ws.Range("M94").Value = -420.75
$ws.Range("N94").Value = -1989.25

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 821.36365
$ws.Range("I22").Value = 754.875
$ws.Range("J22").Value = 998.6667
$ws.Range("K22").Value = 754.875
$ws.Range("L22").Value = 998.6667
$ws.Range("M22").Value = -404.875
$ws.Range("N22").Value = -1698.6667
$ws.Range("H23").Value = 2654.5
$ws.Range("I23").Value = 2654.5
$ws.Range("K23").Value = 2654.5
$ws.Range("M23").Value = -2414.5
$ws.Range("H27").Value = 2654.5
$ws.Range("I27").Value = 2654.5
$ws.Range("K27").Value = 2654.5
$ws.Range("M27").Value = -2462.5
$ws.Range("H31").Value = 93731.266
$ws.Range("I31").Value = 137179.12
$ws.Range("J31").Value = 27111.2
$ws.Range("K31").Value = 137179.12
$ws.Range("L31").Value = 27111.2
$ws.Range("M31").Value = -136884.12
$ws.Range("N31").Value = -27701.2
$ws.Range("H34").Value = 93731.266
$ws.Range("I34").Value = 137179.12
$ws.Range("J34").Value = 27111.2
$ws.Range("K34").Value = 137179.12
$ws.Range("L34").Value = 27111.2
$ws.Range("M34").Value = -136977.12
$ws.Range("N34").Value = -27515.2
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H58").Value = 652633.5
$ws.Range("I58").Value = 1123630.8
$ws.Range("K58").Value = 1123630.8
$ws.Range("M58").Value = -1123427.8
$ws.Range("H136").Value = 652633.5
$ws.Range("I136").Value = 1123630.8
$ws.Range("K136").Value = 3370892.4
$ws.Range("M136").Value = -3368342.4

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 50637372
$ws.Range("I4").Value = 55065296
$ws.Range("K4").Value = 165195888
$ws.Range("M4").Value = -165195776
$ws.Range("H5").Value = 230
$ws.Range("J5").Value = 200
$ws.Range("L5").Value = 600
$ws.Range("N5").Value = -824
$ws.Range("H7").Value = 4435363.5
$ws.Range("I7").Value = 5747143
$ws.Range("J7").Value = 500025
$ws.Range("K7").Value = 17241429
$ws.Range("L7").Value = 1500075
$ws.Range("M7").Value = -17241317
$ws.Range("N7").Value = -1500299
$ws.Range("H135").Value = 230
$ws.Range("J135").Value = 200
$ws.Range("L135").Value = 1800
$ws.Range("N135").Value = -6870

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 1836817.5
$ws.Range("I21").Value = 5009998.5
$ws.Range("K21").Value = 5009998.5
$ws.Range("M21").Value = -5009825.5
$ws.Range("H30").Value = 1836817.5
$ws.Range("I30").Value = 5009998.5
$ws.Range("K30").Value = 5009998.5
$ws.Range("M30").Value = -5009893.5
$ws.Range("H40").Value = 20001
$ws.Range("J40").Value = 20001
$ws.Range("L40").Value = 20001
$ws.Range("N40").Value = -20303
$ws.Range("H44").Value = 3502.2856
$ws.Range("J44").Value = 3502.2856
$ws.Range("L44").Value = 3502.2856
$ws.Range("N44").Value = -4694.2856
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H70").Value = 6456.857
$ws.Range("I70").Value = 6456.857
$ws.Range("K70").Value = 6456.857
$ws.Range("M70").Value = -6186.857
$ws.Range("H73").Value = 6456.857
$ws.Range("I73").Value = 6456.857
$ws.Range("K73").Value = 6456.857
$ws.Range("M73").Value = -5520.857
$ws.Range("H93").Value = 80251
$ws.Range("J93").Value = 80251
$ws.Range("L93").Value = 80251
$ws.Range("N93").Value = -83995
$ws.Range("H132").Value = 19850294
$ws.Range("I132").Value = 25956124
$ws.Range("J132").Value = 6347.8335
$ws.Range("K132").Value = 77868372
$ws.Range("L132").Value = 19043.5005
$ws.Range("M132").Value = -77865842
$ws.Range("N132").Value = -24103.5005

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 663
$ws.Range("I22").Value = 567.4286
$ws.Range("J22").Value = 997.5
$ws.Range("K22").Value = 567.4286
$ws.Range("L22").Value = 997.5
$ws.Range("M22").Value = -272.4286
$ws.Range("N22").Value = -1587.5
$ws.Range("H27").Value = 663
$ws.Range("I27").Value = 567.4286
$ws.Range("J27").Value = 997.5
$ws.Range("K27").Value = 567.4286
$ws.Range("L27").Value = 997.5
$ws.Range("M27").Value = -460.4286
$ws.Range("N27").Value = -1211.5
$ws.Range("H98").Value = 78058.71000000001
$ws.Range("J98").Value = 78058.71000000001
$ws.Range("L98").Value = 78058.71000000001
$ws.Range("N98").Value = -84048.71000000001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 10000
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H132").Value = 7455915
$ws.Range("I132").Value = 9148470
$ws.Range("J132").Value = 8671.6
$ws.Range("K132").Value = 27445410
$ws.Range("L132").Value = 26014.8
$ws.Range("M132").Value = -27442880
$ws.Range("N132").Value = -31074.8
